$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new accommodation row (row 2) beneath the existing header row.
$ws.Range("A2").Value = "Nuevo Apto cerca al aeropuerto con parqueadero."
$ws.Range("B2").Value = "Angee Julieth"
$ws.Range("C2").Value = '$996,889 COP'

# Leave the selection as it was when the workbook was last saved.
$ws.Range("A2:D3").Select()
